# Trade #77 closed at 2026-02-17 08:58:28 - unknown UNKNOWN +0.000%
$wb = $excel.ActiveWorkbook

# --- Summary sheet: refresh aggregate stats ---
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B3").Value = 1200.53   # Current Capital
$wsSummary.Range("B4").Value = 0.54      # Total P&L $
$wsSummary.Range("B6").Value = 77        # Total Trades
$wsSummary.Range("B7").Value = 32        # Winning Trades
$wsSummary.Range("B9").Value = 41.56     # Win Rate %

# --- Strategy Status sheet: refresh MarketMaking row (row 4) ---
$wsStatus = $wb.Worksheets.Item("Strategy Status")
$wsStatus.Range("C4").Value = 100.53     # Capital
$wsStatus.Range("D4").Value = 77         # Trades
$wsStatus.Range("E4").Value = 0.54       # P&L $
$wsStatus.Range("F4").Value = 0.53       # P&L %
$wsStatus.Range("G4").Value = 41.56      # Win Rate %

# --- Append the newly closed trade (#77) as row 78 on both trade log sheets ---
foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Cells.Item(78, 1).Value = 77

    # Force the Date/Time-looking text to stay as plain text instead of being
    # auto-converted into Excel date/time serials.
    $ws.Cells.Item(78, 2).NumberFormat = "@"
    $ws.Cells.Item(78, 2).Value = "2026-02-17"

    $ws.Cells.Item(78, 3).NumberFormat = "@"
    $ws.Cells.Item(78, 3).Value = "08:58:22"

    $ws.Cells.Item(78, 4).Value = "MarketMaking"
    $ws.Cells.Item(78, 5).Value = "DOWN"
    $ws.Cells.Item(78, 6).Value = 0.85
    $ws.Cells.Item(78, 7).Value = 0.86
    $ws.Cells.Item(78, 8).Value = "CLOSED"
    $ws.Cells.Item(78, 9).Value = 1.1765
    $ws.Cells.Item(78, 10).Value = 0.01
    $ws.Cells.Item(78, 11).Value = 100.53
    $ws.Cells.Item(78, 12).Value = 0
    $ws.Cells.Item(78, 13).Value = 0
    $ws.Cells.Item(78, 14).Value = 0.6
    $ws.Cells.Item(78, 15).Value = "Normal spread capture: 19600 bps"
    $ws.Cells.Item(78, 16).Value = "early_exit"
    $ws.Cells.Item(78, 17).Value = 0.1
}
